$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2499.0852
$ws.Range("J17").Value = 2551.0217
$ws.Range("L17").Value = 7653.0651
$ws.Range("N17").Value = -7989.0651

$ws.Range("H113").Value = 2763.8333
$ws.Range("I113").Value = 2000.7778
$ws.Range("J113").Value = 3526.889
$ws.Range("K113").Value = 2000.7778
$ws.Range("L113").Value = 3526.889
$ws.Range("M113").Value = 1253.2222
$ws.Range("N113").Value = -10034.889

$ws.Range("H116").Value = 2300
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -9084

$ws.Range("H137").Value = 2125.5
$ws.Range("I137").Value = 1930.68
$ws.Range("K137").Value = 5792.04
$ws.Range("M137").Value = -3242.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17116.855
$ws.Range("I32").Value = 18615.59
$ws.Range("J32").Value = 5689
$ws.Range("K32").Value = 18615.59
$ws.Range("L32").Value = 5689
$ws.Range("M32").Value = -18328.59
$ws.Range("N32").Value = -6263

$ws.Range("H61").Value = 9732.387000000001
$ws.Range("I61").Value = 5742.5
$ws.Range("J61").Value = 16986.727
$ws.Range("K61").Value = 5742.5
$ws.Range("L61").Value = 16986.727
$ws.Range("M61").Value = -5530.5
$ws.Range("N61").Value = -17410.727

$ws.Range("H74").Value = 1961.1389
$ws.Range("I74").Value = 1938.3438
$ws.Range("J74").Value = 2143.5
$ws.Range("K74").Value = 1938.3438
$ws.Range("L74").Value = 2143.5
$ws.Range("M74").Value = -1064.3438
$ws.Range("N74").Value = -3891.5

$ws.Range("H77").Value = 1961.1389
$ws.Range("I77").Value = 1938.3438
$ws.Range("J77").Value = 2143.5
$ws.Range("K77").Value = 9691.719000000001
$ws.Range("L77").Value = 10717.5
$ws.Range("M77").Value = -5323.719000000001
$ws.Range("N77").Value = -19453.5

$ws.Range("H132").Value = 2734.4
$ws.Range("I132").Value = 2242.6667
$ws.Range("J132").Value = 3188.3076
$ws.Range("K132").Value = 6728.000100000001
$ws.Range("L132").Value = 9564.9228
$ws.Range("M132").Value = -4198.000100000001
$ws.Range("N132").Value = -14624.9228

$ws.Range("H136").Value = 9732.387000000001
$ws.Range("I136").Value = 5742.5
$ws.Range("J136").Value = 16986.727
$ws.Range("K136").Value = 17227.5
$ws.Range("L136").Value = 50960.181
$ws.Range("M136").Value = -14677.5
$ws.Range("N136").Value = -56060.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 896.8889
$ws.Range("I107").Value = 881.35
$ws.Range("J107").Value = 941.2857
$ws.Range("K107").Value = 881.35
$ws.Range("L107").Value = 941.2857
$ws.Range("M107").Value = 1038.65
$ws.Range("N107").Value = -4781.2857

$ws.Range("H134").Value = 44135.457
$ws.Range("I134").Value = 2627.3125
$ws.Range("J134").Value = 127151.75
$ws.Range("K134").Value = 7881.9375
$ws.Range("L134").Value = 381455.25
$ws.Range("M134").Value = -5346.9375
$ws.Range("N134").Value = -386525.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2418.3057
$ws.Range("I31").Value = 1574.8334
$ws.Range("K31").Value = 1574.8334
$ws.Range("M31").Value = -1279.8334

$ws.Range("H34").Value = 2418.3057
$ws.Range("I34").Value = 1574.8334
$ws.Range("K34").Value = 1574.8334
$ws.Range("M34").Value = -1372.8334

$ws.Range("H58").Value = 2397222.5
$ws.Range("I58").Value = 3638270.2
$ws.Range("J58").Value = 10592.462
$ws.Range("K58").Value = 3638270.2
$ws.Range("L58").Value = 10592.462
$ws.Range("M58").Value = -3638067.2
$ws.Range("N58").Value = -10998.462

$ws.Range("H132").Value = 2533.0476
$ws.Range("I132").Value = 2305.6775
$ws.Range("J132").Value = 3173.818
$ws.Range("K132").Value = 6917.032499999999
$ws.Range("L132").Value = 9521.454000000002
$ws.Range("M132").Value = -4387.032499999999
$ws.Range("N132").Value = -14581.454

$ws.Range("H134").Value = 2279.3333
$ws.Range("I134").Value = 2139.7812
$ws.Range("J134").Value = 2725.9
$ws.Range("K134").Value = 6419.3436
$ws.Range("L134").Value = 8177.700000000001
$ws.Range("M134").Value = -3884.3436
$ws.Range("N134").Value = -13247.7

$ws.Range("H136").Value = 2397222.5
$ws.Range("I136").Value = 3638270.2
$ws.Range("J136").Value = 10592.462
$ws.Range("K136").Value = 10914810.6
$ws.Range("L136").Value = 31777.386
$ws.Range("M136").Value = -10912260.6
$ws.Range("N136").Value = -36877.386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79799
$ws.Range("J37").Value = 79799
$ws.Range("L37").Value = 239397
$ws.Range("N37").Value = -239621

$ws.Range("H60").Value = 423.33334
$ws.Range("I60").Value = 316.9
$ws.Range("J60").Value = 955.5
$ws.Range("K60").Value = 950.6999999999999
$ws.Range("L60").Value = 2866.5
$ws.Range("M60").Value = -699.6999999999999
$ws.Range("N60").Value = -3368.5

$ws.Range("H98").Value = 362.5263
$ws.Range("I98").Value = 310.4
$ws.Range("J98").Value = 558
$ws.Range("K98").Value = 931.1999999999999
$ws.Range("L98").Value = 1674
$ws.Range("M98").Value = 566.8000000000001
$ws.Range("N98").Value = -4670

$ws.Range("H110").Value = 3839.1177
$ws.Range("I110").Value = 1600
$ws.Range("J110").Value = 3979.0625
$ws.Range("K110").Value = 4800
$ws.Range("L110").Value = 11937.1875
$ws.Range("M110").Value = -710
$ws.Range("N110").Value = -20117.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13767.818
$ws.Range("I132").Value = 6026.8
$ws.Range("J132").Value = 20218.666
$ws.Range("K132").Value = 18080.4
$ws.Range("L132").Value = 60655.99800000001
$ws.Range("M132").Value = -15550.4
$ws.Range("N132").Value = -65715.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3240.5134
$ws.Range("I132").Value = 2788.138
$ws.Range("J132").Value = 4880.375
$ws.Range("K132").Value = 8364.414000000001
$ws.Range("L132").Value = 14641.125
$ws.Range("M132").Value = -5834.414000000001
$ws.Range("N132").Value = -19701.125

$ws.Range("H136").Value = 4156.393
$ws.Range("I136").Value = 2663.3
$ws.Range("J136").Value = 5879.1924
$ws.Range("K136").Value = 7989.900000000001
$ws.Range("L136").Value = 17637.5772
$ws.Range("M136").Value = -5439.900000000001
$ws.Range("N136").Value = -22737.5772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1756.56
$ws.Range("I132").Value = 1010.7692
$ws.Range("K132").Value = 3032.3076
$ws.Range("M132").Value = -502.3076000000001

$ws.Range("H136").Value = 2571.432
$ws.Range("I136").Value = 1638.1
$ws.Range("J136").Value = 4571.4287
$ws.Range("K136").Value = 4914.299999999999
$ws.Range("L136").Value = 13714.2861
$ws.Range("M136").Value = -2364.299999999999
$ws.Range("N136").Value = -18814.2861
